$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New user rows appended to the table (rows 38 & 39)
$ws.Range("A38").Value = "سارا اکبری"
$ws.Range("B38").Value = "آموزگار"
$ws.Range("C38").Value = "sara"

$ws.Range("A39").Value = "پرهام یزدانی"
$ws.Range("B39").Value = "والد"
$ws.Range("C39").Value = "parham"

# Update the view: scrolled down a bit further and a new active selection
$ws.Range("E38").Select()
